# Updated cryptos list (price + 1h volume change) for each coin row.
# Numeric-looking "Price" values are entered with a leading apostrophe so
# Excel keeps them as literal text (e.g. "157.30" instead of 157.3),
# matching values that already contain more than one '.' and therefore
# stay text automatically (e.g. "42.460.96").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.460.96'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '2.288.33'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'157.30"
$ws.Range('E5').Value = '  +15,621.00%  '
$ws.Range('D6').Value = "'307.65"
$ws.Range('E6').Value = '  +1.48%  '
$ws.Range('D7').Value = "'95.90"
$ws.Range('E7').Value = '  +5.39%  '
$ws.Range('E8').Value = '  +0.65%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('E10').Value = '  +3.49%  '
$ws.Range('D11').Value = "'36.06"
$ws.Range('E11').Value = '  +12.53%  '
$ws.Range('E12').Value = '  +1.31%  '
$ws.Range('E13').Value = '  -1.92%  '
$ws.Range('E14').Value = '  +2.80%  '
$ws.Range('D15').Value = '2.642.78'
$ws.Range('E15').Value = '  +1.24%  '
$ws.Range('D16').Value = "'14.52"
$ws.Range('E16').Value = '  +2.84%  '
$ws.Range('D17').Value = '2.283.37'
$ws.Range('E17').Value = '  +0.71%  '
$ws.Range('D18').Value = "'0.804"
$ws.Range('E18').Value = '  +6.22%  '
$ws.Range('D19').Value = '42.390.32'
$ws.Range('E19').Value = '  +2.02%  '
$ws.Range('D20').Value = "'12.66"
$ws.Range('E20').Value = '  +2.49%  '
$ws.Range('E21').Value = '  +2.17%  '
$ws.Range('D22').Value = "'6.01"
$ws.Range('E22').Value = '  +2.24%  '
$ws.Range('D23').Value = "'68.08"
$ws.Range('E23').Value = '  +2.47%  '
$ws.Range('D24').Value = "'243.08"
$ws.Range('E24').Value = '  +1.29%  '
$ws.Range('E25').Value = '  +1.41%  '
$ws.Range('E26').Value = '  +2.17%  '
$ws.Range('E27').Value = '  -0.19%  '
$ws.Range('D28').Value = "'24.06"
$ws.Range('E28').Value = '  +0.26%  '
$ws.Range('D29').Value = "'35.98"
$ws.Range('E29').Value = '  +5.32%  '
$ws.Range('D30').Value = "'9.58"
$ws.Range('E30').Value = '  +1.32%  '
$ws.Range('E31').Value = '  +1.39%  '
$ws.Range('D32').Value = "'160.78"
$ws.Range('E32').Value = '  +0.42%  '
$ws.Range('E33').Value = '  +4.00%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  +2.16%  '
$ws.Range('E36').Value = '  +3.36%  '
$ws.Range('E37').Value = '  +4.96%  '
$ws.Range('D38').Value = "'17.29"
$ws.Range('E38').Value = '  +4.75%  '
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('E40').Value = '  +4.33%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').Value = "'4.17"
$ws.Range('E42').Value = '  +7.20%  '
$ws.Range('D43').Value = '2.009.40'
$ws.Range('E43').Value = '  -2.41%  '
$ws.Range('E44').Value = '  +12.09%  '
$ws.Range('E45').Value = '  -0.44%  '
$ws.Range('E46').Value = '  +3.30%  '
$ws.Range('E47').Value = '  +5.88%  '
$ws.Range('E48').Value = '  +0.57%  '
$ws.Range('D49').Value = "'53.86"
$ws.Range('E49').Value = '  +5.02%  '
$ws.Range('E50').Value = '  +2.99%  '
$ws.Range('D51').Value = "'72.93"
$ws.Range('E51').Value = '  +0.87%  '
